$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.845.88'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '3.369.32'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '662.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.44'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.431'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -4.17%  '
$ws.Range("D11").Value = '3.367.90'
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  -1.28%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = '97.793.48'
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.16'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.93%  '
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("D17").Value = '3.995.10'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.99%  '
$ws.Range("D19").Value = '3.369.08'
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.549'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +9.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.90'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '516.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.02%  '
$ws.Range("E25").Value = '  -1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +12.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '94.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.90%  '
$ws.Range("D29").Value = '3.541.90'
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.148'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.24%  '
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("E33").Value = '  -6.36%  '
$ws.Range("E34").Value = '  +14.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.577'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '28.98'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.05%  '
$ws.Range("E39").Value = '  +6.84%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '524.26'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.153'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.47%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0441'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.02%  '
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.864'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("E46").Value = '  +10.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.98%  '
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.09%  '
